$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the last row (currently row 22, "Cleanup code")
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

# Row 22: Revisit end game
$ws.Range("A22").Value = "Revisit end game, more elegant solution"
$ws.Range("B22").Value = 0.0625
$ws.Range("B22").Style = $ws.Range("B2").Style
$ws.Range("C22").Value = "DONE"
$ws.Range("C22").Style = $ws.Range("C2").Style
$ws.Range("D22").Value = "This cost me a lot of time due to testing"

# Row 23: Inventory (panel)
$ws.Range("A23").Value = "Inventory ( panel )"
$ws.Range("B23").Value = 0.0069444444444444441
$ws.Range("B23").Style = $ws.Range("B2").Style
$ws.Range("C23").Value = "DONE"
$ws.Range("C23").Style = $ws.Range("C2").Style
$ws.Range("D23").Value = "Seperated character panel and inventory panel ( put actual items on the inventory panel )"

$ws.Range("D23").Select()
